# Auto-generated Excel COM-interop script
# Applies updated market price / profit figures to the Leviathan_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1645.6154
$ws.Range("I28").Value = 1454
$ws.Range("J28").Value = 2699.5
$ws.Range("K28").Value = 1454
$ws.Range("L28").Value = 2699.5
$ws.Range("M28").Value = -969
$ws.Range("N28").Value = -3669.5
$ws.Range("H33").Value = 641.7778
$ws.Range("I33").Value = 641.7778
$ws.Range("K33").Value = 641.7778
$ws.Range("M33").Value = -412.7778
$ws.Range("H86").Value = 1887.6154
$ws.Range("I86").Value = 1939
$ws.Range("J86").Value = 1748.1428
$ws.Range("K86").Value = 1939
$ws.Range("L86").Value = 1748.1428
$ws.Range("M86").Value = -816
$ws.Range("N86").Value = -3994.1428
$ws.Range("H89").Value = 1887.6154
$ws.Range("I89").Value = 1939
$ws.Range("J89").Value = 1748.1428
$ws.Range("K89").Value = 9695
$ws.Range("L89").Value = 8740.714
$ws.Range("M89").Value = -4079
$ws.Range("N89").Value = -19972.714
$ws.Range("H92").Value = 548.5
$ws.Range("I92").Value = 412.94116
$ws.Range("K92").Value = 412.94116
$ws.Range("M92").Value = 835.0588399999999
$ws.Range("H94").Value = 16670401
$ws.Range("I94").Value = 25002600
$ws.Range("K94").Value = 25002600
$ws.Range("M94").Value = -25002149
$ws.Range("H98").Value = 1046.8334
$ws.Range("I98").Value = 1046.8334
$ws.Range("K98").Value = 1046.8334
$ws.Range("M98").Value = 451.1666
$ws.Range("H107").Value = 852.4375
$ws.Range("I107").Value = 885.93335
$ws.Range("K107").Value = 885.93335
$ws.Range("M107").Value = 1034.06665
$ws.Range("H112").Value = 1471.6111
$ws.Range("J112").Value = 1754.5454
$ws.Range("L112").Value = 5263.6362
$ws.Range("N112").Value = -7479.6362
$ws.Range("H122").Value = 1046.8334
$ws.Range("I122").Value = 1046.8334
$ws.Range("K122").Value = 3140.5002
$ws.Range("M122").Value = -690.5001999999999
$ws.Range("H132").Value = 2398.8147
$ws.Range("I132").Value = 1790.76
$ws.Range("K132").Value = 5372.28
$ws.Range("M132").Value = -2842.28
$ws.Range("H135").Value = 1521.2778
$ws.Range("I135").Value = 1384.0667
$ws.Range("K135").Value = 12456.6003
$ws.Range("M135").Value = -9921.6003
$ws.Range("H138").Value = 2490.7334
$ws.Range("I138").Value = 1659.6875
$ws.Range("J138").Value = 3440.5
$ws.Range("K138").Value = 4979.0625
$ws.Range("L138").Value = 10321.5
$ws.Range("M138").Value = 160.9375
$ws.Range("N138").Value = -20601.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 223111.33
$ws.Range("I4").Value = 250875.25
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 250875.25
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -250759.25
$ws.Range("N4").Value = -1232
$ws.Range("H34").Value = 24500
$ws.Range("J34").Value = 24500
$ws.Range("L34").Value = 24500
$ws.Range("N34").Value = -25042

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2359
$ws.Range("J7").Value = 4997.5
$ws.Range("L7").Value = 4997.5
$ws.Range("N7").Value = -5223.5
$ws.Range("H86").Value = 2327.6
$ws.Range("I86").Value = 929.3333
$ws.Range("K86").Value = 929.3333
$ws.Range("M86").Value = 193.6667
$ws.Range("H89").Value = 2327.6
$ws.Range("I89").Value = 929.3333
$ws.Range("K89").Value = 4646.6665
$ws.Range("M89").Value = 969.3334999999997
$ws.Range("H105").Value = 4213.846
$ws.Range("I105").Value = 4963.2
$ws.Range("K105").Value = 4963.2
$ws.Range("M105").Value = -3216.2
$ws.Range("H107").Value = 550
$ws.Range("I107").Value = 550
$ws.Range("K107").Value = 550
$ws.Range("M107").Value = 1370

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 36666668
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 678.05884
$ws.Range("I5").Value = 425.33334
$ws.Range("J5").Value = 962.375
$ws.Range("K5").Value = 1276.00002
$ws.Range("L5").Value = 2887.125
$ws.Range("M5").Value = -1164.00002
$ws.Range("N5").Value = -3111.125
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H113").Value = 1063.8182
$ws.Range("J113").Value = 866.8889
$ws.Range("L113").Value = 2600.6667
$ws.Range("N113").Value = -6940.6667
$ws.Range("H129").Value = 68864.836
$ws.Range("J129").Value = 2566.0908
$ws.Range("L129").Value = 7698.2724
$ws.Range("N129").Value = -17698.2724
$ws.Range("H135").Value = 678.05884
$ws.Range("I135").Value = 425.33334
$ws.Range("J135").Value = 962.375
$ws.Range("K135").Value = 3828.00006
$ws.Range("L135").Value = 8661.375
$ws.Range("M135").Value = -1293.00006
$ws.Range("N135").Value = -13731.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 285.42856
$ws.Range("I2").Value = 294
$ws.Range("K2").Value = 294
$ws.Range("M2").Value = -181
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -388
$ws.Range("N4").ClearContents()
$ws.Range("H11").Value = 7709267.5
$ws.Range("I11").Value = 4324857.5
$ws.Range("J11").Value = 10868050
$ws.Range("K11").Value = 4324857.5
$ws.Range("L11").Value = 10868050
$ws.Range("M11").Value = -4324718.5
$ws.Range("N11").Value = -10868328
$ws.Range("H80").Value = 4178.2856
$ws.Range("I80").Value = 3130.6
$ws.Range("K80").Value = 3130.6
$ws.Range("M80").Value = -2132.6
$ws.Range("H83").Value = 4178.2856
$ws.Range("I83").Value = 3130.6
$ws.Range("K83").Value = 15653
$ws.Range("M83").Value = -10661
$ws.Range("H107").Value = 62501890
$ws.Range("I107").Value = 1150.5
$ws.Range("J107").Value = 125002620
$ws.Range("K107").Value = 1150.5
$ws.Range("L107").Value = 125002620
$ws.Range("M107").Value = 769.5
$ws.Range("N107").Value = -125006460
$ws.Range("H113").Value = 3247.8125
$ws.Range("J113").Value = 5954
$ws.Range("L113").Value = 5954
$ws.Range("N113").Value = -10294

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 755
$ws.Range("I22").Value = 741.5714
$ws.Range("J22").Value = 778.5
$ws.Range("K22").Value = 741.5714
$ws.Range("L22").Value = 778.5
$ws.Range("M22").Value = -446.5714
$ws.Range("N22").Value = -1368.5
$ws.Range("H27").Value = 755
$ws.Range("I27").Value = 741.5714
$ws.Range("J27").Value = 778.5
$ws.Range("K27").Value = 741.5714
$ws.Range("L27").Value = 778.5
$ws.Range("M27").Value = -634.5714
$ws.Range("N27").Value = -992.5
$ws.Range("H40").Value = 9997.714
$ws.Range("I40").Value = 10331.5
$ws.Range("K40").Value = 10331.5
$ws.Range("M40").Value = -10195.5
$ws.Range("H46").Value = 1468
$ws.Range("I46").Value = 1594
$ws.Range("K46").Value = 1594
$ws.Range("M46").Value = -1406
$ws.Range("H61").Value = 54312.26
$ws.Range("I61").Value = 64090.562
$ws.Range("J61").Value = 2161.3333
$ws.Range("K61").Value = 64090.562
$ws.Range("L61").Value = 2161.3333
$ws.Range("M61").Value = -63888.562
$ws.Range("N61").Value = -2565.3333
$ws.Range("H82").Value = 3864.75
$ws.Range("I82").Value = 3855.7334
$ws.Range("K82").Value = 3855.7334
$ws.Range("M82").Value = -3494.7334
$ws.Range("H85").Value = 3864.75
$ws.Range("I85").Value = 3855.7334
$ws.Range("K85").Value = 3855.7334
$ws.Range("M85").Value = -2607.7334
$ws.Range("H93").Value = 1774.2727
$ws.Range("I93").Value = 1774.2727
$ws.Range("K93").Value = 1774.2727
$ws.Range("M93").Value = -526.2727
$ws.Range("H113").Value = 54312.26
$ws.Range("I113").Value = 64090.562
$ws.Range("J113").Value = 2161.3333
$ws.Range("K113").Value = 64090.562
$ws.Range("L113").Value = 2161.3333
$ws.Range("M113").Value = -61920.562
$ws.Range("N113").Value = -6501.3333
$ws.Range("H122").Value = 14553.333
$ws.Range("I122").Value = 27247
$ws.Range("K122").Value = 81741
$ws.Range("M122").Value = -79291
$ws.Range("H132").Value = 2997
$ws.Range("I132").Value = 2640
$ws.Range("J132").Value = 7995
$ws.Range("K132").Value = 7920
$ws.Range("L132").Value = 23985
$ws.Range("M132").Value = -5390
$ws.Range("N132").Value = -29045

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6600.1113
$ws.Range("I81").Value = 6950.125
$ws.Range("J81").Value = 3800
$ws.Range("K81").Value = 13900.25
$ws.Range("L81").Value = 7600
$ws.Range("M81").Value = -12839.25
$ws.Range("N81").Value = -9722
$ws.Range("H84").Value = 6600.1113
$ws.Range("I84").Value = 6950.125
$ws.Range("J84").Value = 3800
$ws.Range("K84").Value = 69501.25
$ws.Range("L84").Value = 38000
$ws.Range("M84").Value = -64197.25
$ws.Range("N84").Value = -48608
$ws.Range("H122").Value = 2953.3845
$ws.Range("I122").Value = 3045.3333
$ws.Range("K122").Value = 9135.999899999999
$ws.Range("M122").Value = -6685.999899999999
$ws.Range("H132").Value = 40001
$ws.Range("I132").Value = 40001
$ws.Range("K132").Value = 120003
$ws.Range("M132").Value = -117473
